$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from the row above into the new row's date cell
$ws.Range("A7").Copy($ws.Range("A8"))

$ws.Range("A8").Value = 42916
$ws.Range("B8").Value = "Kelly Marinduque"
$ws.Range("C8").Value = "JObstreet"
$ws.Range("D8").Value = "Fiber Technician"
$ws.Range("E8").Value = 943254
